$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 ("Time Period"): roll the end-of-range month from 2023-11 to 2023-12
# (Nuclear / column N keeps its existing range, unchanged this update)
$ws.Cells.Item(6, 2).Value  = "2008-12:2023-12"   # B6
$ws.Cells.Item(6, 3).Value  = "2008-12:2023-12"   # C6
$ws.Cells.Item(6, 4).Value  = "2010-09:2023-12"   # D6
$ws.Cells.Item(6, 5).Value  = "2008-12:2023-12"   # E6
$ws.Cells.Item(6, 6).Value  = "2009-12:2023-12"   # F6
$ws.Cells.Item(6, 7).Value  = "2017-12:2023-12"   # G6
$ws.Cells.Item(6, 8).Value  = "2017-12:2023-12"   # H6
$ws.Cells.Item(6, 9).Value  = "2012-12:2023-12"   # I6
$ws.Cells.Item(6, 10).Value = "2012-12:2023-12"   # J6
$ws.Cells.Item(6, 11).Value = "2012-12:2023-12"   # K6
$ws.Cells.Item(6, 12).Value = "2012-12:2023-12"   # L6
$ws.Cells.Item(6, 13).Value = "2012-12:2023-12"   # M6
$ws.Cells.Item(6, 14).Value = "2020-06:2023-11"   # N6 (unchanged)
$ws.Cells.Item(6, 15).Value = "2017-12:2023-12"   # O6
$ws.Cells.Item(6, 16).Value = "2017-12:2023-12"   # P6
$ws.Cells.Item(6, 17).Value = "2018-02:2023-12"   # Q6

# --- Row 8 ("Update"): refresh the update-stamp dates
# (leading apostrophe forces these to stay plain text, like the source
# cells, instead of being auto-converted to Excel date serials)
$ws.Cells.Item(8, 2).Value  = "'2024-02-05"   # B8
$ws.Cells.Item(8, 3).Value  = "'2024-02-05"   # C8
$ws.Cells.Item(8, 4).Value  = "'2024-02-05"   # D8
$ws.Cells.Item(8, 5).Value  = "'2024-02-05"   # E8
$ws.Cells.Item(8, 6).Value  = "'2024-02-05"   # F8
$ws.Cells.Item(8, 7).Value  = "'2024-02-05"   # G8
$ws.Cells.Item(8, 8).Value  = "'2024-02-05"   # H8
$ws.Cells.Item(8, 9).Value  = "'2024-01-29"   # I8
$ws.Cells.Item(8, 10).Value = "'2024-01-29"   # J8
$ws.Cells.Item(8, 11).Value = "'2024-01-29"   # K8
$ws.Cells.Item(8, 12).Value = "'2024-01-29"   # L8
$ws.Cells.Item(8, 13).Value = "'2024-01-29"   # M8
$ws.Cells.Item(8, 14).Value = "'2023-12-29"   # N8 (unchanged)
$ws.Cells.Item(8, 15).Value = "'2024-02-05"   # O8
$ws.Cells.Item(8, 16).Value = "'2024-02-05"   # P8
$ws.Cells.Item(8, 17).Value = "'2024-02-05"   # Q8

# --- Append a new data row (172) for period 2023-12-31, copying the
# number/date formatting down from the last existing row (171)
$ws.Range("A171:Q171").Copy()
$ws.Range("A172").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(172, 1).Value  = 45291   # date serial for 2023-12-31
$ws.Cells.Item(172, 2).Value  = 1034
$ws.Cells.Item(172, 3).Value  = 6567
$ws.Cells.Item(172, 4).Value  = 139
$ws.Cells.Item(172, 5).Value  = 7566
$ws.Cells.Item(172, 6).Value  = 21602
$ws.Cells.Item(172, 7).Value  = 4774
$ws.Cells.Item(172, 8).Value  = 1012
$ws.Cells.Item(172, 9).Value  = 291965
$ws.Cells.Item(172, 10).Value = 42154
$ws.Cells.Item(172, 11).Value = 139032
$ws.Cells.Item(172, 12).Value = 5691
$ws.Cells.Item(172, 13).Value = 44134
$ws.Cells.Item(172, 14).Value = 0
$ws.Cells.Item(172, 15).Value = 116493
$ws.Cells.Item(172, 16).Value = 12562
$ws.Cells.Item(172, 17).Value = 60949
